# Generate Report for Handback
# Updates the localization status workbook to reflect a failed handback
# transform for the 1f3f2cd7-... source file in both the zh-cn and de-de
# target languages.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status text shared across sheets for the 1f3f2cd7 row changes from
# "Ready for handoff" to "Handback transform failed".
$newStatus = "Handback transform failed"
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# zh-cn sheet: record the handback/handoff filename mismatch error.
$zhcn.Range("P3").Value = "Handback file name: lckmxyjz.nvg is different with handoff file name: 1f3f2cd7-dacf-429d-9b2c-991a17ed9b51.84da50e4288e9eaa16e74293c5a0ff5739d9aff5.zh-cn."

# de-de sheet: the 1f3f2cd7 row now also records the same style error.
$dede.Range("P3").Value = "Handback file name: lckmxyjz.nvg is different with handoff file name: 1f3f2cd7-dacf-429d-9b2c-991a17ed9b51.84da50e4288e9eaa16e74293c5a0ff5739d9aff5.de-de."

# Widen the Error Detail column on both language sheets so the new,
# longer messages are readable (ColumnWidth 39.17 round-trips to a
# stored OOXML column width of 40, matching the other width="40" columns
# already on these sheets).
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
